$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "false start" data rows (original rows 2 and 3), shifting
# the remaining rows up.
$ws.Rows("2:3").Delete()

# Restore the selection state that Excel would have after this edit
# (both remaining data rows, selected as whole rows).
$ws.Range("A2:XFD3").Select()
